# Update "想去人数" (attendance count) figures on the 展览 and 全部类型 sheets
# to match the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6976
$ws1.Range("F12").Value = 197
$ws1.Range("F15").Value = 1816
$ws1.Range("F17").Value = 3588
$ws1.Range("F19").Value = 243
$ws1.Range("F23").Value = 2206
$ws1.Range("F24").Value = 3
$ws1.Range("F30").Value = 16
$ws1.Range("F32").Value = 130
$ws1.Range("F33").Value = 64

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6976
$ws4.Range("F13").Value = 197
$ws4.Range("F16").Value = 1816
$ws4.Range("F18").Value = 3588
$ws4.Range("F20").Value = 243
$ws4.Range("F24").Value = 2206
$ws4.Range("F25").Value = 3
$ws4.Range("F31").Value = 16
$ws4.Range("F33").Value = 131
$ws4.Range("F34").Value = 64
